$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear row 3 (old record) entirely first
$ws.Range("A3:D3").Clear()

# Remove row4 leftover styled cells entirely
$ws.Range("A4:D4").Clear()

# Update row 2 values
$ws.Range("D2").Value = "OEK0628"
$ws.Range("A2").Value = "4575000"

# C2: change number format from Text (@) to Number (0.00), keep value
$ws.Range("C2").NumberFormat = "0.00"

# A3 literal value 0, with the leftover text-format style
$ws.Range("A3").Value = 0
$ws.Range("A3").NumberFormat = "@"

# Set selection to A2
$ws.Range("A2").Select()
